$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.932.10"
$ws.Range("E2").Value = "  +1.59%  "
$ws.Range("D3").Value = "1.641.44"
$ws.Range("E3").Value = "  +0.76%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.96"
$ws.Range("E5").Value = "  +0.73%  "
$ws.Range("E6").Value = "  +0.55%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.47"
$ws.Range("E8").Value = "  +2.19%  "
$ws.Range("E9").Value = "  -1.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0614"
$ws.Range("E10").Value = "  +0.55%  "
$ws.Range("E11").Value = "  +2.52%  "
$ws.Range("E12").Value = "  +0.51%  "
$ws.Range("D13").Value = "1.642.87"
$ws.Range("E13").Value = "  +0.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.09"
$ws.Range("E14").Value = "  +1.19%  "
$ws.Range("E15").Value = "  +2.55%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.59"
$ws.Range("E16").Value = "  +0.95%  "
$ws.Range("D17").Value = "27.931.00"
$ws.Range("E17").Value = "  +1.50%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "232.98"
$ws.Range("E18").Value = "  +1.93%  "
$ws.Range("D19").Value = "0.0₃0722"
$ws.Range("E19").Value = "  +0.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.60"
$ws.Range("E20").Value = "  +0.68%  "
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.49"
$ws.Range("E22").Value = "  -2.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.38"
$ws.Range("E23").Value = "  +0.47%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.10"
$ws.Range("E24").Value = "  -1.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.85"
$ws.Range("E25").Value = "  +2.56%  "
$ws.Range("E26").Value = "  +0.55%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.71"
$ws.Range("E27").Value = "  +0.93%  "
$ws.Range("E28").Value = "  +0.39%  "
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("E30").Value = "  +0.83%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0485"
$ws.Range("E31").Value = "  +0.77%  "
$ws.Range("E32").Value = "  +2.82%  "
$ws.Range("E33").Value = "  +0.65%  "
$ws.Range("D34").Value = "1.408.67"
$ws.Range("E34").Value = "  -3.73%  "
$ws.Range("E35").Value = "  +2.57%  "
$ws.Range("E36").Value = "  +1.56%  "
$ws.Range("E37").Value = "  +1.73%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.882"
$ws.Range("E38").Value = "  +0.77%  "
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.560"
$ws.Range("E39").Value = "  +0.47%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.925"
$ws.Range("E40").Value = "  +0.94%  "
$ws.Range("E41").Value = "  +1.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "67.53"
$ws.Range("E43").Value = "  -0.59%  "
$ws.Range("E44").Value = "  +6.35%  "
$ws.Range("E45").Value = "  +2.84%  "
$ws.Range("E46").Value = "  +0.09%  "
$ws.Range("D47").Value = "1.781.01"
$ws.Range("E47").Value = "  +0.59%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.83"
$ws.Range("E48").Value = "  +0.56%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.100"
$ws.Range("E49").Value = "  +0.78%  "
$ws.Range("E50").Value = "  +0.44%  "
$ws.Range("E51").Value = "  -0.05%  "
